# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 9952
$wsExhibit.Range("F10").Value = 4987
$wsExhibit.Range("F17").Value = 292
$wsExhibit.Range("F19").Value = 113
$wsExhibit.Range("F21").Value = 1494

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 9952
$wsAll.Range("F11").Value = 4987
$wsAll.Range("F18").Value = 292
$wsAll.Range("F20").Value = 113
$wsAll.Range("F22").Value = 1494
